$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Female")
$ws.Range("H3").Value = 38529
$ws.Range("O3").Value = 8131.2
$ws.Range("Q3").Value = 26.7
$ws.Range("T3").Value = "8131.2 (±444.2)"
$ws.Range("U3").Value = "26.7% (±1.8%)"
$ws.Range("W3").Value = 1992.1
$ws.Range("Y3").Value = "1992.1(±108.8)"
$ws.Range("H6").Value = 31218
$ws.Range("O6").Value = 5437
$ws.Range("T6").Value = "5437.0 (±474.0)"
$ws.Range("W6").Value = 1892
$ws.Range("Y6").Value = "1892.0(±165.0)"
$ws.Range("H8").Value = 4639
$ws.Range("O8").Value = 455
$ws.Range("Q8").Value = 10.9
$ws.Range("R8").Value = 2.3
$ws.Range("T8").Value = "455.0 (±87.0)"
$ws.Range("U8").Value = "10.9% (±2.3%)"
$ws.Range("W8").Value = 818.3
$ws.Range("Y8").Value = "818.3(±156.5)"
$ws.Range("H10").Value = 201761
$ws.Range("O10").Value = 27019.6
$ws.Range("Q10").Value = 15.5
$ws.Range("R10").Value = 1.1
$ws.Range("T10").Value = "27019.6 (±1669.3)"
$ws.Range("U10").Value = "15.5% (±1.1%)"
$ws.Range("W10").Value = 1041.1
$ws.Range("Y10").Value = "1041.1(±64.3)"
$ws.Range("H14").Value = 7890
$ws.Range("O14").Value = 781.8
$ws.Range("Q14").Value = 11
$ws.Range("T14").Value = "781.8 (±116.6)"
$ws.Range("U14").Value = "11.0% (±1.8%)"
$ws.Range("W14").Value = 969.3
$ws.Range("X14").Value = 144.5
$ws.Range("Y14").Value = "969.3(±144.5)"
$ws.Range("H18").Value = 1011
$ws.Range("O18").Value = 177.4
$ws.Range("Q18").Value = 21.3
$ws.Range("R18").Value = 2.1
$ws.Range("T18").Value = "177.4 (±14.2)"
$ws.Range("U18").Value = "21.3% (±2.1%)"
$ws.Range("W18").Value = 1377.2
$ws.Range("Y18").Value = "1377.2(±110.3)"
$ws.Range("H19").Value = 1433
$ws.Range("O19").Value = 191
$ws.Range("Q19").Value = 15.4
$ws.Range("T19").Value = "191.0 (±46.5)"
$ws.Range("U19").Value = "15.4% (±4.2%)"
$ws.Range("W19").Value = 1569.9
$ws.Range("X19").Value = 382.3
$ws.Range("Y19").Value = "1569.9(±382.3)"
$ws.Range("D20").Value = 39143
$ws.Range("E20").Value = 38421
$ws.Range("F20").Value = 39124
$ws.Range("G20").Value = 39494
$ws.Range("H20").Value = 45695
$ws.Range("I20").Value = 39738.79999999999
$ws.Range("J20").Value = 357.9
$ws.Range("L20").Value = 313.7
$ws.Range("M20").Value = 39425.09999999999
$ws.Range("N20").Value = 40052.49999999999
$ws.Range("O20").Value = 5956.2
$ws.Range("P20").Value = 313.7
$ws.Range("Q20").Value = 15
$ws.Range("S20").Value = "39738.8 (±313.7)"
$ws.Range("T20").Value = "5956.2 (±313.7)"
$ws.Range("U20").Value = "15.0% (±0.9%)"
$ws.Range("W20").Value = 1213.4
$ws.Range("Y20").Value = "1213.4(±64.0)"
$ws.Range("H21").Value = 11076
$ws.Range("O21").Value = -270.2
$ws.Range("T21").Value = "-270.2 (±177.4)"
$ws.Range("W21").Value = -195.1
$ws.Range("Y21").Value = "-195.1(±128.1)"
$ws.Range("H24").Value = 63785
$ws.Range("O24").Value = 9495.4
$ws.Range("Q24").Value = 17.5
$ws.Range("T24").Value = "9495.4 (±1037.8)"
$ws.Range("U24").Value = "17.5% (±2.2%)"
$ws.Range("W24").Value = 1609.7
$ws.Range("Y24").Value = "1609.7(±176.0)"
$ws.Range("H26").Value = 12831
$ws.Range("O26").Value = 1498.2
$ws.Range("T26").Value = "1498.2 (±199.3)"
$ws.Range("W26").Value = 1200.5
$ws.Range("Y26").Value = "1200.5(±159.7)"
$ws.Range("H27").Value = 7327
$ws.Range("O27").Value = 1841.6
$ws.Range("Q27").Value = 33.6
$ws.Range("T27").Value = "1841.6 (±110.7)"
$ws.Range("U27").Value = "33.6% (±2.7%)"
$ws.Range("W27").Value = 2468.3
$ws.Range("X27").Value = 148.3
$ws.Range("Y27").Value = "2468.3(±148.3)"
$ws.Range("H28").Value = 151378
$ws.Range("O28").Value = 30137
$ws.Range("Q28").Value = 24.9
$ws.Range("R28").Value = 0.8
$ws.Range("T28").Value = "30137.0 (±731.0)"
$ws.Range("U28").Value = "24.9% (±0.8%)"
$ws.Range("W28").Value = 1662.8
$ws.Range("Y28").Value = "1662.8(±40.4)"

$ws = $wb.Worksheets.Item("Male")
$ws.Range("H3").Value = 26690
$ws.Range("O3").Value = 6140.4
$ws.Range("Q3").Value = 29.9
$ws.Range("T3").Value = "6140.4 (±289.3)"
$ws.Range("U3").Value = "29.9% (±1.8%)"
$ws.Range("W3").Value = 2571.5
$ws.Range("Y3").Value = "2571.5(±121.1)"
$ws.Range("H6").Value = 19770
$ws.Range("O6").Value = 4591.8
$ws.Range("Q6").Value = 30.3
$ws.Range("R6").Value = 2.2
$ws.Range("T6").Value = "4591.8 (±258.7)"
$ws.Range("U6").Value = "30.3% (±2.2%)"
$ws.Range("W6").Value = 3155
$ws.Range("Y6").Value = "3155.0(±177.7)"
$ws.Range("H8").Value = 2058
$ws.Range("O8").Value = 270.8
$ws.Range("Q8").Value = 15.2
$ws.Range("T8").Value = "270.8 (±64.3)"
$ws.Range("U8").Value = "15.2% (±4.0%)"
$ws.Range("W8").Value = 1416.5
$ws.Range("X8").Value = 336.4
$ws.Range("Y8").Value = "1416.5(±336.4)"
$ws.Range("H10").Value = 140555
$ws.Range("O10").Value = 22004.2
$ws.Range("Q10").Value = 18.6
$ws.Range("R10").Value = 1.3
$ws.Range("T10").Value = "22004.2 (±1234.7)"
$ws.Range("U10").Value = "18.6% (±1.3%)"
$ws.Range("W10").Value = 1531.2
$ws.Range("Y10").Value = "1531.2(±85.9)"
$ws.Range("H14").Value = 3302
$ws.Range("O14").Value = 416.2
$ws.Range("Q14").Value = 14.4
$ws.Range("R14").Value = 3.1
$ws.Range("T14").Value = "416.2 (±81.6)"
$ws.Range("U14").Value = "14.4% (±3.1%)"
$ws.Range("W14").Value = 1549.9
$ws.Range("X14").Value = 303.8
$ws.Range("Y14").Value = "1549.9(±303.8)"
$ws.Range("H19").Value = 1078
$ws.Range("O19").Value = 219.4
$ws.Range("Q19").Value = 25.6
$ws.Range("T19").Value = "219.4 (±35.0)"
$ws.Range("U19").Value = "25.6% (±5.0%)"
$ws.Range("W19").Value = 2717.7
$ws.Range("X19").Value = 433.5
$ws.Range("Y19").Value = "2717.7(±433.5)"
$ws.Range("D20").Value = 26821
$ws.Range("E20").Value = 26470
$ws.Range("F20").Value = 27648
$ws.Range("G20").Value = 28749
$ws.Range("H20").Value = 35182
$ws.Range("I20").Value = 27720.99999999999
$ws.Range("M20").Value = 26910.69999999999
$ws.Range("N20").Value = 28531.29999999999
$ws.Range("O20").Value = 7461
$ws.Range("R20").Value = 3.6
$ws.Range("S20").Value = "27721.0 (±810.3)"
$ws.Range("T20").Value = "7461.0 (±810.3)"
$ws.Range("U20").Value = "26.9% (±3.6%)"
$ws.Range("W20").Value = 2422.7
$ws.Range("Y20").Value = "2422.7(±263.1)"
$ws.Range("H21").Value = 7988
$ws.Range("O21").Value = 141.2
$ws.Range("Q21").Value = 1.8
$ws.Range("R21").Value = 1.2
$ws.Range("T21").Value = "141.2 (±97.3)"
$ws.Range("U21").Value = "1.8% (±1.2%)"
$ws.Range("W21").Value = 161.4
$ws.Range("Y21").Value = "161.4(±111.2)"
$ws.Range("H24").Value = 42812
$ws.Range("O24").Value = 8665.200000000001
$ws.Range("Q24").Value = 25.4
$ws.Range("R24").Value = 1.9
$ws.Range("T24").Value = "8665.2 (±516.5)"
$ws.Range("U24").Value = "25.4% (±1.9%)"
$ws.Range("W24").Value = 2737.6
$ws.Range("Y24").Value = "2737.6(±163.2)"
$ws.Range("H27").Value = 4333
$ws.Range("O27").Value = 1232
$ws.Range("Q27").Value = 39.7
$ws.Range("T27").Value = "1232.0 (±96.2)"
$ws.Range("U27").Value = "39.7% (±4.2%)"
$ws.Range("W27").Value = 3382.6
$ws.Range("X27").Value = 264.1
$ws.Range("Y27").Value = "3382.6(±264.1)"
$ws.Range("H28").Value = 110136
$ws.Range("O28").Value = 22667
$ws.Range("Q28").Value = 25.9
$ws.Range("T28").Value = "22667.0 (±651.6)"
$ws.Range("U28").Value = "25.9% (±0.9%)"
$ws.Range("W28").Value = 2121.4
$ws.Range("Y28").Value = "2121.4(±61.0)"

$ws = $wb.Worksheets.Item("Total")
$ws.Range("H3").Value = 65219
$ws.Range("O3").Value = 14271.6
$ws.Range("Q3").Value = 28
$ws.Range("R3").Value = 1.6
$ws.Range("T3").Value = "14271.6 (±653.5)"
$ws.Range("U3").Value = "28.0% (±1.6%)"
$ws.Range("W3").Value = 2205.9
$ws.Range("Y3").Value = "2205.9(±101.0)"
$ws.Range("H6").Value = 50988
$ws.Range("O6").Value = 10028.8
$ws.Range("T6").Value = "10028.8 (±672.3)"
$ws.Range("W6").Value = 2316.6
$ws.Range("Y6").Value = "2316.6(±155.3)"
$ws.Range("H8").Value = 6697
$ws.Range("O8").Value = 725.8
$ws.Range("Q8").Value = 12.2
$ws.Range("R8").Value = 2.8
$ws.Range("T8").Value = "725.8 (±148.8)"
$ws.Range("U8").Value = "12.2% (±2.8%)"
$ws.Range("W8").Value = 971.4
$ws.Range("Y8").Value = "971.4(±199.2)"
$ws.Range("H10").Value = 342316
$ws.Range("O10").Value = 49023.8
$ws.Range("T10").Value = "49023.8 (±2877.7)"
$ws.Range("W10").Value = 1215.7
$ws.Range("X10").Value = 71.40000000000001
$ws.Range("Y10").Value = "1215.7(±71.4)"
$ws.Range("H14").Value = 11192
$ws.Range("O14").Value = 1198
$ws.Range("Q14").Value = 12
$ws.Range("T14").Value = "1198.0 (±175.6)"
$ws.Range("U14").Value = "12.0% (±1.9%)"
$ws.Range("W14").Value = 1114.3
$ws.Range("X14").Value = 163.3
$ws.Range("Y14").Value = "1114.3(±163.3)"
$ws.Range("H18").Value = 1744
$ws.Range("O18").Value = 319.4
$ws.Range("Q18").Value = 22.4
$ws.Range("T18").Value = "319.4 (±37.1)"
$ws.Range("U18").Value = "22.4% (±3.1%)"
$ws.Range("W18").Value = 1537.6
$ws.Range("Y18").Value = "1537.6(±178.6)"
$ws.Range("H19").Value = 2511
$ws.Range("O19").Value = 410.4
$ws.Range("Q19").Value = 19.5
$ws.Range("R19").Value = 3.9
$ws.Range("T19").Value = "410.4 (±72.2)"
$ws.Range("U19").Value = "19.5% (±3.9%)"
$ws.Range("W19").Value = 2027.8
$ws.Range("Y19").Value = "2027.8(±356.7)"
$ws.Range("D20").Value = 65964
$ws.Range("E20").Value = 64891
$ws.Range("F20").Value = 66772
$ws.Range("G20").Value = 68242
$ws.Range("H20").Value = 80877
$ws.Range("I20").Value = 67459.59999999999
$ws.Range("J20").Value = 1233.5
$ws.Range("L20").Value = 1081.2
$ws.Range("M20").Value = 66378.39999999999
$ws.Range("N20").Value = 68540.79999999999
$ws.Range("O20").Value = 13417.4
$ws.Range("P20").Value = 1081.2
$ws.Range("Q20").Value = 19.9
$ws.Range("S20").Value = "67459.6 (±1081.2)"
$ws.Range("T20").Value = "13417.4 (±1081.2)"
$ws.Range("U20").Value = "19.9% (±1.9%)"
$ws.Range("W20").Value = 1679.7
$ws.Range("X20").Value = 135.3
$ws.Range("Y20").Value = "1679.7(±135.3)"
$ws.Range("H21").Value = 19064
$ws.Range("O21").Value = -129
$ws.Range("Q21").Value = -0.7
$ws.Range("T21").Value = "-129.0 (±231.8)"
$ws.Range("U21").Value = "-0.7% (±1.2%)"
$ws.Range("W21").Value = -57.1
$ws.Range("X21").Value = 102.6
$ws.Range("Y21").Value = "-57.1(±102.6)"
$ws.Range("H24").Value = 106597
$ws.Range("O24").Value = 18160.6
$ws.Range("T24").Value = "18160.6 (±1494.8)"
$ws.Range("W24").Value = 2003.6
$ws.Range("Y24").Value = "2003.6(±164.9)"
$ws.Range("H26").Value = 19898
$ws.Range("O26").Value = 2734.2
$ws.Range("T26").Value = "2734.2 (±293.1)"
$ws.Range("W26").Value = 1523.1
$ws.Range("X26").Value = 163.3
$ws.Range("Y26").Value = "1523.1(±163.3)"
$ws.Range("H27").Value = 11660
$ws.Range("O27").Value = 3073.6
$ws.Range("Q27").Value = 35.8
$ws.Range("T27").Value = "3073.6 (±177.1)"
$ws.Range("U27").Value = "35.8% (±2.7%)"
$ws.Range("W27").Value = 2768.2
$ws.Range("Y27").Value = "2768.2(±159.5)"
$ws.Range("H28").Value = 261514
$ws.Range("O28").Value = 52804
$ws.Range("Q28").Value = 25.3
$ws.Range("T28").Value = "52804.0 (±1249.5)"
$ws.Range("U28").Value = "25.3% (±0.7%)"
$ws.Range("W28").Value = 1832.9
$ws.Range("Y28").Value = "1832.9(±43.4)"
